# Update column E (rows 2-23) from 50 to 70
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2:E23").Value = 70

# Update the active cell/selection to L27 (single cell), matching the diff
$ws.Range("L27").Select()
